# This script updates the cryptocurrency price/volume table (columns B-E, rows 2-51)
# to match the refreshed data feed values in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, even when the text looks
# like a number (e.g. "1.003" or "1.000"), preserving the original default style.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.585.20"
$ws.Range("E2").Value = "  +2.32%  "

# Row 3
Set-TextValue "D3" "1.995.22"
$ws.Range("E3").Value = "  +6.16%  "

# Row 4
Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextValue "D5" "326.46"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
Set-TextValue "D7" "0.4692"
$ws.Range("E7").Value = "  +1.73%  "

# Row 8
Set-TextValue "D8" "0.3956"
$ws.Range("E8").Value = "  +2.23%  "

# Row 9
Set-TextValue "D9" "46.57"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
Set-TextValue "D10" "0.07961"
$ws.Range("E10").Value = "  +1.42%  "

# Row 11
Set-TextValue "D11" "1.003"
$ws.Range("E11").Value = "  +1.74%  "

# Row 12
Set-TextValue "D12" "22.95"
$ws.Range("E12").Value = "  +5.46%  "

# Row 13
Set-TextValue "D13" "1.979.75"
$ws.Range("E13").Value = "  +6.40%  "

# Row 14
Set-TextValue "D14" "7.281"
$ws.Range("E14").Value = "  +4.22%  "

# Row 15
Set-TextValue "D15" "5.892"
$ws.Range("E15").Value = "  +4.33%  "

# Row 16
Set-TextValue "D16" "0.07138"
$ws.Range("E16").Value = "  +2.31%  "

# Row 17
Set-TextValue "D17" "88.97"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18
Set-TextValue "D18" "1.003"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
Set-TextValue "D19" "0.000009998"
$ws.Range("E19").Value = "  +0.44%  "

# Row 20
Set-TextValue "D20" "17.46"
$ws.Range("E20").Value = "  +3.04%  "

# Row 21
Set-TextValue "D21" "1.000"

# Row 22
Set-TextValue "D22" "29.591.86"
$ws.Range("E22").Value = "  +2.32%  "

# Row 23
Set-TextValue "D23" "5.562"
$ws.Range("E23").Value = "  +6.08%  "

# Row 24
Set-TextValue "D24" "11.29"
$ws.Range("E24").Value = "  +3.10%  "

# Row 25
Set-TextValue "D25" "2.103"

# Row 26
Set-TextValue "D26" "158.01"
$ws.Range("E26").Value = "  +1.01%  "

# Row 27
Set-TextValue "D27" "19.72"
$ws.Range("E27").Value = "  +1.95%  "

# Row 28
Set-TextValue "D28" "6.013"
$ws.Range("E28").Value = "  +0.40%  "

# Row 29
Set-TextValue "D29" "120.32"
$ws.Range("E29").Value = "  +2.40%  "

# Row 30
Set-TextValue "D30" "1.957"
$ws.Range("E30").Value = "  +1.96%  "

# Row 31
Set-TextValue "D31" "0.09459"
$ws.Range("E31").Value = "  +0.95%  "

# Row 32
Set-TextValue "D32" "0.9147"
$ws.Range("E32").Value = "  +1.55%  "

# Row 33
Set-TextValue "D33" "1.355"
$ws.Range("E33").Value = "  +2.87%  "

# Row 34
Set-TextValue "D34" "5.270"
$ws.Range("E34").Value = "  +0.11%  "

# Row 35
Set-TextValue "D35" "3.182"
$ws.Range("E35").Value = "  -2.21%  "

# Row 36
Set-TextValue "D36" "0.000003502"
$ws.Range("E36").Value = "  +100.21%  "

# Row 37
Set-TextValue "D37" "0.05867"
$ws.Range("E37").Value = "  +2.14%  "

# Row 38
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
Set-TextValue "D39" "0.02125"
$ws.Range("E39").Value = "  +2.54%  "

# Row 40
Set-TextValue "D40" "7.943"
$ws.Range("E40").Value = "  +3.99%  "

# Row 41
Set-TextValue "D41" "0.5783"
$ws.Range("E41").Value = "  +2.40%  "

# Row 42
Set-TextValue "D42" "0.1833"
$ws.Range("E42").Value = "  +3.83%  "

# Row 43
Set-TextValue "D43" "9.865"
$ws.Range("E43").Value = "  +1.70%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "12.07"
$ws.Range("E44").Value = "  +1.07%  "

# Row 45
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D45" "2.758"
$ws.Range("E45").Value = "  +8.39%  "

# Row 46
Set-TextValue "D46" "0.5394"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
Set-TextValue "D47" "2.195"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
Set-TextValue "D48" "1.880"
$ws.Range("E48").Value = "  +1.99%  "

# Row 49
Set-TextValue "D49" "0.06950"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
Set-TextValue "D50" "114.42"
$ws.Range("E50").Value = "  +1.42%  "

# Row 51
Set-TextValue "D51" "0.3084"
$ws.Range("E51").Value = "  +7.30%  "
